# Generate Report for handback
# Update the "Correspond Handoff Datetime" (D2) and "Correspond Handback
# DateTime" (G2) entries for the first file row on both the zh-cn and
# de-de report sheets, reflecting the newly generated handback report
# timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-21 02:29:05"
$wsZhCn.Range("G2").Value = "2016-01-21 02:29:52"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-21 02:29:17"
$wsDeDe.Range("G2").Value = "2016-01-21 02:30:14"
